$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the per-record data among rows 2,3,4,5,6,9 (rows 7 and 8
# are untouched). Concretely: new row 2 <- old row 9, new row 3 <- old row 5,
# new row 4 <- old row 3, new row 5 <- old row 2, new row 6 <- old row 4,
# new row 9 <- old row 6. Capture every source row's values BEFORE writing
# anything, then write the captured snapshots into their destination rows.

function Get-RowSnapshot($r) {
    $snap = @{
        A  = $ws.Range("A$r").Value2
        B  = $ws.Range("B$r").Value2
        D  = $ws.Range("D$r").Value2
        E  = $ws.Range("E$r").Value2
        F  = $ws.Range("F$r").Value2
        G  = $ws.Range("G$r").Value2
        H  = $ws.Range("H$r").Value2
        I  = $ws.Range("I$r").Text
        J  = $ws.Range("J$r").Value2
        K  = $ws.Range("K$r").Value2
        Q  = $ws.Range("Q$r").Value2
        R  = $ws.Range("R$r").Value2
        AC = $ws.Range("AC$r").Value2
    }
    return $snap
}

$row2 = Get-RowSnapshot 2
$row3 = Get-RowSnapshot 3
$row4 = Get-RowSnapshot 4
$row5 = Get-RowSnapshot 5
$row6 = Get-RowSnapshot 6
$row9 = Get-RowSnapshot 9

function Set-RowSnapshot($r, $snap) {
    $ws.Range("A$r").Value = $snap.A
    $ws.Range("B$r").Value = $snap.B
    $ws.Range("D$r").Value = $snap.D
    $ws.Range("E$r").Value = $snap.E
    $ws.Range("F$r").Value = $snap.F
    $ws.Range("G$r").Value = $snap.G
    $ws.Range("H$r").Value = $snap.H

    # Column I holds digit-only labels ("30", "15", "2", ...) that must stay
    # TEXT, not be coerced to numbers - force text format before assigning.
    $ws.Range("I$r").NumberFormat = "@"
    if ($snap.I -eq $null -or $snap.I -eq "") {
        $ws.Range("I$r").ClearContents()
        $ws.Range("I$r").NumberFormat = "@"
    } else {
        $ws.Range("I$r").Value = [string]$snap.I
    }

    if ($snap.J -eq $null -or $snap.J -eq "") {
        $ws.Range("J$r").ClearContents()
    } else {
        $ws.Range("J$r").Value = $snap.J
    }

    if ($snap.K -eq $null -or $snap.K -eq "") {
        $ws.Range("K$r").ClearContents()
    } else {
        $ws.Range("K$r").Value = $snap.K
    }

    $ws.Range("Q$r").Value = $snap.Q
    $ws.Range("R$r").Value = $snap.R

    if ($snap.AC -eq $null -or $snap.AC -eq "") {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $snap.AC
    }
}

Set-RowSnapshot 2 $row9
Set-RowSnapshot 3 $row5
Set-RowSnapshot 4 $row3
Set-RowSnapshot 5 $row2
Set-RowSnapshot 6 $row4
Set-RowSnapshot 9 $row6
